$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts old E:J to F:K)
$ws.Range("E1").EntireColumn.Insert()

# Header
$ws.Range("E1").Value = "frequency"
$ws.Range("E1:E3").NumberFormat = $ws.Range("F1").NumberFormat

# Data values
$ws.Range("E2").Value = "MONTHLY"
$ws.Range("E3").Value = "QUARTERLY"

$ws.Range("E3").Select()
